$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range('D2:E2')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '245.85'
$vals[0,1] = '0.75%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D3:E3')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '29.21'
$vals[0,1] = '6.34%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D4:E4')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '5.179'
$vals[0,1] = '1.06%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D5:E5')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.05732'
$vals[0,1] = '0.97%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D6:E6')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '6.569'
$vals[0,1] = '0.32%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('B7:E7')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,4
$vals[0,0] = 'MXToken'
$vals[0,1] = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$vals[0,2] = '0.8606'
$vals[0,3] = '5.03%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('B8:E8')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,4
$vals[0,0] = 'FTXToken'
$vals[0,1] = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$vals[0,2] = '0.8618'
$vals[0,3] = '0.42%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('B9:E9')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,4
$vals[0,0] = 'WazirX'
$vals[0,1] = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$vals[0,2] = '0.1365'
$vals[0,3] = '1.80%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('B10:E10')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,4
$vals[0,0] = 'MandalaExchangeToken'
$vals[0,1] = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$vals[0,2] = '0.07069'
$vals[0,3] = '1.70%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('B11:E11')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,4
$vals[0,0] = 'BitrueCoin'
$vals[0,1] = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$vals[0,2] = '0.03063'
$vals[0,3] = '6.53%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('B12:E12')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,4
$vals[0,0] = 'BitMartToken'
$vals[0,1] = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$vals[0,2] = '0.09376'
$vals[0,3] = '-0.21%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('B13:E13')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,4
$vals[0,0] = 'BitForexToken'
$vals[0,1] = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$vals[0,2] = '0.001540'
$vals[0,3] = '1.25%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('B14:E14')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,4
$vals[0,0] = 'One'
$vals[0,1] = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$vals[0,2] = '0.0006005'
$vals[0,3] = '0.04%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('B15:E15')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,4
$vals[0,0] = 'TigerCash'
$vals[0,1] = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$vals[0,2] = '0.005999'
$vals[0,3] = '-3.49%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('B16:E16')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,4
$vals[0,0] = 'UpBots'
$vals[0,1] = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$vals[0,2] = '0.007489'
$vals[0,3] = '5,223.82%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('B17:E17')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,4
$vals[0,0] = 'LEO'
$vals[0,1] = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$vals[0,2] = '3.489'
$vals[0,3] = '-0.76%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('B18:E18')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,4
$vals[0,0] = 'GateToken'
$vals[0,1] = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$vals[0,2] = '3.108'
$vals[0,3] = '3.28%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D19:E19')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '2.185'
$vals[0,1] = '0.72%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('E20')
$rng.NumberFormat = "@"
$rng.Value = '1.60%'
$rng.ClearFormats()

$rng = $ws.Range('E21')
$rng.NumberFormat = "@"
$rng.Value = '2.73%'
$rng.ClearFormats()

$rng = $ws.Range('E22')
$rng.NumberFormat = "@"
$rng.Value = '-1.54%'
$rng.ClearFormats()

$rng = $ws.Range('D23:E23')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '3.479'
$vals[0,1] = '-2.83%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D24:E24')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.04149'
$vals[0,1] = '1.08%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('E25')
$rng.NumberFormat = "@"
$rng.Value = '0.46%'
$rng.ClearFormats()

$rng = $ws.Range('D26:E26')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.001228'
$vals[0,1] = '1.05%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D27:E27')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.004988'
$vals[0,1] = '11.61%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D28:E28')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.0001211'
$vals[0,1] = '2.65%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D40:E40')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.03739'
$vals[0,1] = '0.58%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D41:E41')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.005795'
$vals[0,1] = '-2.42%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D42:E42')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.1070'
$vals[0,1] = '1.37%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D43:E43')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.002462'
$vals[0,1] = '6.79%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D44:E44')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.008613'
$vals[0,1] = '-4.81%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D45:E45')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.00005282'
$vals[0,1] = '3.38%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D46:E46')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.00000000751'
$vals[0,1] = '0.11%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D47')
$rng.NumberFormat = "@"
$rng.Value = '0.05705'
$rng.ClearFormats()

$rng = $ws.Range('D48:E48')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.002255'
$vals[0,1] = '-10.97%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D49:E49')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.00002102'
$vals[0,1] = '0.11%'
$rng.Value = $vals
$rng.ClearFormats()

$rng = $ws.Range('D50:E50')
$rng.NumberFormat = "@"
$vals = New-Object "string[,]" 1,2
$vals[0,0] = '0.0002002'
$vals[0,1] = '0.11%'
$rng.Value = $vals
$rng.ClearFormats()
